# Re-applies the row permutation + content fix described by the commit diff.
# Rows 2,3,5-18 have their entire row content (except the row number itself)
# replaced with the content that, before the edit, lived in a different row
# (row 4 is untouched by the edit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: set to the content that was previously in row 13
$ws.Range("A2").Value = 111396324
$ws.Range("B2").Value = 96348
$ws.Range("C2").Value = '''Ovaliderad'
$ws.Range("D2").Value = '''VU'
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = '''Knärot'
$ws.Range("G2").Value = '''Goodyera repens'
$ws.Range("H2").Value = '''(L.) R. Br.'
$ws.Range("P2").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q2").Value = 625335.6676841485
$ws.Range("R2").Value = 7209609.168182318
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = '''Västerbotten'
$ws.Range("U2").Value = '''Storuman'
$ws.Range("V2").Value = '''Lycksele lappmark'
$ws.Range("W2").Value = '''Stensele'
$ws.Range("Y2").Value = '''2023-08-08'
$ws.Range("Z2").Value = '''00:00'
$ws.Range("AA2").Value = '''2023-08-08'
$ws.Range("AB2").Value = '''00:00'
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AW2").Value = '''Isak Vahlström'
$ws.Range("AX2").Value = '''Via Isak Vahlström'

# Row 3: set to the content that was previously in row 12
$ws.Range("A3").Value = 111396308
$ws.Range("B3").Value = 56398
$ws.Range("C3").Value = '''Ovaliderad'
$ws.Range("D3").Value = '''NT'
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = '''Tretåig hackspett'
$ws.Range("G3").Value = '''Picoides tridactylus'
$ws.Range("H3").Value = '''(Linnaeus, 1758)'
$ws.Range("M3").Value = '''äldre spår'
$ws.Range("P3").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q3").Value = 625151.1577179903
$ws.Range("R3").Value = 7209567.512248591
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = '''Västerbotten'
$ws.Range("U3").Value = '''Storuman'
$ws.Range("V3").Value = '''Lycksele lappmark'
$ws.Range("W3").Value = '''Stensele'
$ws.Range("Y3").Value = '''2023-08-08'
$ws.Range("Z3").Value = '''00:00'
$ws.Range("AA3").Value = '''2023-08-08'
$ws.Range("AB3").Value = '''00:00'
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = '''Isak Vahlström'
$ws.Range("AX3").Value = '''Via Isak Vahlström'

# Row 5: set to the content that was previously in row 3
$ws.Range("A5").Value = 111396313
$ws.Range("B5").Value = 96348
$ws.Range("C5").Value = '''Ovaliderad'
$ws.Range("D5").Value = '''VU'
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = '''Knärot'
$ws.Range("G5").Value = '''Goodyera repens'
$ws.Range("H5").Value = '''(L.) R. Br.'
$ws.Range("P5").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q5").Value = 625231.5510770321
$ws.Range("R5").Value = 7209481.895207534
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = '''Västerbotten'
$ws.Range("U5").Value = '''Storuman'
$ws.Range("V5").Value = '''Lycksele lappmark'
$ws.Range("W5").Value = '''Stensele'
$ws.Range("Y5").Value = '''2023-08-08'
$ws.Range("Z5").Value = '''00:00'
$ws.Range("AA5").Value = '''2023-08-08'
$ws.Range("AB5").Value = '''00:00'
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = '''Isak Vahlström'
$ws.Range("AX5").Value = '''Via Isak Vahlström'

# Row 6: set to the content that was previously in row 8
$ws.Range("A6").Value = 111396325
$ws.Range("B6").Value = 96348
$ws.Range("C6").Value = '''Ovaliderad'
$ws.Range("D6").Value = '''VU'
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = '''Knärot'
$ws.Range("G6").Value = '''Goodyera repens'
$ws.Range("H6").Value = '''(L.) R. Br.'
$ws.Range("P6").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q6").Value = 625389.9085714296
$ws.Range("R6").Value = 7209580.514361567
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = '''Västerbotten'
$ws.Range("U6").Value = '''Storuman'
$ws.Range("V6").Value = '''Lycksele lappmark'
$ws.Range("W6").Value = '''Stensele'
$ws.Range("Y6").Value = '''2023-08-08'
$ws.Range("Z6").Value = '''00:00'
$ws.Range("AA6").Value = '''2023-08-08'
$ws.Range("AB6").Value = '''00:00'
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = '''Isak Vahlström'
$ws.Range("AX6").Value = '''Via Isak Vahlström'

# Row 7: set to the content that was previously in row 16
$ws.Range("A7").Value = 111396310
$ws.Range("B7").Value = 96348
$ws.Range("C7").Value = '''Ovaliderad'
$ws.Range("D7").Value = '''VU'
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = '''Knärot'
$ws.Range("G7").Value = '''Goodyera repens'
$ws.Range("H7").Value = '''(L.) R. Br.'
$ws.Range("P7").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q7").Value = 625289.0018867656
$ws.Range("R7").Value = 7209518.212698339
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = '''Västerbotten'
$ws.Range("U7").Value = '''Storuman'
$ws.Range("V7").Value = '''Lycksele lappmark'
$ws.Range("W7").Value = '''Stensele'
$ws.Range("Y7").Value = '''2023-08-08'
$ws.Range("Z7").Value = '''00:00'
$ws.Range("AA7").Value = '''2023-08-08'
$ws.Range("AB7").Value = '''00:00'
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = '''Isak Vahlström'
$ws.Range("AX7").Value = '''Via Isak Vahlström'

# Row 8: set to the content that was previously in row 17
$ws.Range("A8").Value = 111396321
$ws.Range("B8").Value = 96348
$ws.Range("C8").Value = '''Ovaliderad'
$ws.Range("D8").Value = '''VU'
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = '''Knärot'
$ws.Range("G8").Value = '''Goodyera repens'
$ws.Range("H8").Value = '''(L.) R. Br.'
$ws.Range("P8").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q8").Value = 625240.2002264742
$ws.Range("R8").Value = 7209649.650274927
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = '''Västerbotten'
$ws.Range("U8").Value = '''Storuman'
$ws.Range("V8").Value = '''Lycksele lappmark'
$ws.Range("W8").Value = '''Stensele'
$ws.Range("Y8").Value = '''2023-08-08'
$ws.Range("Z8").Value = '''00:00'
$ws.Range("AA8").Value = '''2023-08-08'
$ws.Range("AB8").Value = '''00:00'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = '''Isak Vahlström'
$ws.Range("AX8").Value = '''Via Isak Vahlström'

# Row 9: set to the content that was previously in row 6
$ws.Range("A9").Value = 111396316
$ws.Range("B9").Value = 96348
$ws.Range("C9").Value = '''Ovaliderad'
$ws.Range("D9").Value = '''VU'
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = '''Knärot'
$ws.Range("G9").Value = '''Goodyera repens'
$ws.Range("H9").Value = '''(L.) R. Br.'
$ws.Range("P9").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q9").Value = 625153.7279882778
$ws.Range("R9").Value = 7209526.513740451
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = '''Västerbotten'
$ws.Range("U9").Value = '''Storuman'
$ws.Range("V9").Value = '''Lycksele lappmark'
$ws.Range("W9").Value = '''Stensele'
$ws.Range("Y9").Value = '''2023-08-08'
$ws.Range("Z9").Value = '''00:00'
$ws.Range("AA9").Value = '''2023-08-08'
$ws.Range("AB9").Value = '''00:00'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AW9").Value = '''Isak Vahlström'
$ws.Range("AX9").Value = '''Via Isak Vahlström'

# Row 10: set to the content that was previously in row 15
$ws.Range("A10").Value = 111396317
$ws.Range("B10").Value = 96348
$ws.Range("C10").Value = '''Ovaliderad'
$ws.Range("D10").Value = '''VU'
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = '''Knärot'
$ws.Range("G10").Value = '''Goodyera repens'
$ws.Range("H10").Value = '''(L.) R. Br.'
$ws.Range("P10").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q10").Value = 625153.5624699651
$ws.Range("R10").Value = 7209550.662191558
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = '''Västerbotten'
$ws.Range("U10").Value = '''Storuman'
$ws.Range("V10").Value = '''Lycksele lappmark'
$ws.Range("W10").Value = '''Stensele'
$ws.Range("Y10").Value = '''2023-08-08'
$ws.Range("Z10").Value = '''00:00'
$ws.Range("AA10").Value = '''2023-08-08'
$ws.Range("AB10").Value = '''00:00'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AW10").Value = '''Isak Vahlström'
$ws.Range("AX10").Value = '''Via Isak Vahlström'

# Row 11: set to the content that was previously in row 7
$ws.Range("A11").Value = 111396326
$ws.Range("B11").Value = 96348
$ws.Range("C11").Value = '''Ovaliderad'
$ws.Range("D11").Value = '''VU'
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = '''Knärot'
$ws.Range("G11").Value = '''Goodyera repens'
$ws.Range("H11").Value = '''(L.) R. Br.'
$ws.Range("P11").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q11").Value = 625397.1584455093
$ws.Range("R11").Value = 7209589.718691397
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = '''Västerbotten'
$ws.Range("U11").Value = '''Storuman'
$ws.Range("V11").Value = '''Lycksele lappmark'
$ws.Range("W11").Value = '''Stensele'
$ws.Range("Y11").Value = '''2023-08-08'
$ws.Range("Z11").Value = '''00:00'
$ws.Range("AA11").Value = '''2023-08-08'
$ws.Range("AB11").Value = '''00:00'
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AW11").Value = '''Isak Vahlström'
$ws.Range("AX11").Value = '''Via Isak Vahlström'

# Row 12: set to the content that was previously in row 5
$ws.Range("A12").Value = 111396314
$ws.Range("B12").Value = 96348
$ws.Range("C12").Value = '''Ovaliderad'
$ws.Range("D12").Value = '''VU'
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = '''Knärot'
$ws.Range("G12").Value = '''Goodyera repens'
$ws.Range("H12").Value = '''(L.) R. Br.'
$ws.Range("P12").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q12").Value = 625202.8383709632
$ws.Range("R12").Value = 7209539.171001118
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = '''Västerbotten'
$ws.Range("U12").Value = '''Storuman'
$ws.Range("V12").Value = '''Lycksele lappmark'
$ws.Range("W12").Value = '''Stensele'
$ws.Range("Y12").Value = '''2023-08-08'
$ws.Range("Z12").Value = '''00:00'
$ws.Range("AA12").Value = '''2023-08-08'
$ws.Range("AB12").Value = '''00:00'
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AW12").Value = '''Isak Vahlström'
$ws.Range("AX12").Value = '''Via Isak Vahlström'
$ws.Range("M12").ClearContents()

# Row 13: set to the content that was previously in row 2
$ws.Range("A13").Value = 111396322
$ws.Range("B13").Value = 96348
$ws.Range("C13").Value = '''Ovaliderad'
$ws.Range("D13").Value = '''VU'
$ws.Range("E13").Value = 220787
$ws.Range("F13").Value = '''Knärot'
$ws.Range("G13").Value = '''Goodyera repens'
$ws.Range("H13").Value = '''(L.) R. Br.'
$ws.Range("P13").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q13").Value = 625269.4478252844
$ws.Range("R13").Value = 7209630.115016816
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = '''Västerbotten'
$ws.Range("U13").Value = '''Storuman'
$ws.Range("V13").Value = '''Lycksele lappmark'
$ws.Range("W13").Value = '''Stensele'
$ws.Range("Y13").Value = '''2023-08-08'
$ws.Range("Z13").Value = '''00:00'
$ws.Range("AA13").Value = '''2023-08-08'
$ws.Range("AB13").Value = '''00:00'
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AW13").Value = '''Isak Vahlström'
$ws.Range("AX13").Value = '''Via Isak Vahlström'

# Row 14: set to the content that was previously in row 10
$ws.Range("A14").Value = 111396312
$ws.Range("B14").Value = 96348
$ws.Range("C14").Value = '''Ovaliderad'
$ws.Range("D14").Value = '''VU'
$ws.Range("E14").Value = 220787
$ws.Range("F14").Value = '''Knärot'
$ws.Range("G14").Value = '''Goodyera repens'
$ws.Range("H14").Value = '''(L.) R. Br.'
$ws.Range("P14").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q14").Value = 625242.7087276473
$ws.Range("R14").Value = 7209468.80281719
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = '''Västerbotten'
$ws.Range("U14").Value = '''Storuman'
$ws.Range("V14").Value = '''Lycksele lappmark'
$ws.Range("W14").Value = '''Stensele'
$ws.Range("Y14").Value = '''2023-08-08'
$ws.Range("Z14").Value = '''00:00'
$ws.Range("AA14").Value = '''2023-08-08'
$ws.Range("AB14").Value = '''00:00'
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AW14").Value = '''Isak Vahlström'
$ws.Range("AX14").Value = '''Via Isak Vahlström'

# Row 15: set to the content that was previously in row 11
$ws.Range("A15").Value = 111396323
$ws.Range("B15").Value = 96348
$ws.Range("C15").Value = '''Ovaliderad'
$ws.Range("D15").Value = '''VU'
$ws.Range("E15").Value = 220787
$ws.Range("F15").Value = '''Knärot'
$ws.Range("G15").Value = '''Goodyera repens'
$ws.Range("H15").Value = '''(L.) R. Br.'
$ws.Range("P15").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q15").Value = 625301.6605433678
$ws.Range("R15").Value = 7209610.70454926
$ws.Range("S15").Value = 10
$ws.Range("T15").Value = '''Västerbotten'
$ws.Range("U15").Value = '''Storuman'
$ws.Range("V15").Value = '''Lycksele lappmark'
$ws.Range("W15").Value = '''Stensele'
$ws.Range("Y15").Value = '''2023-08-08'
$ws.Range("Z15").Value = '''00:00'
$ws.Range("AA15").Value = '''2023-08-08'
$ws.Range("AB15").Value = '''00:00'
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AG15").Value = $false
$ws.Range("AW15").Value = '''Isak Vahlström'
$ws.Range("AX15").Value = '''Via Isak Vahlström'

# Row 16: set to the content that was previously in row 14
$ws.Range("A16").Value = 111396311
$ws.Range("B16").Value = 96348
$ws.Range("C16").Value = '''Ovaliderad'
$ws.Range("D16").Value = '''VU'
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = '''Knärot'
$ws.Range("G16").Value = '''Goodyera repens'
$ws.Range("H16").Value = '''(L.) R. Br.'
$ws.Range("P16").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q16").Value = 625271.0561409625
$ws.Range("R16").Value = 7209511.101565193
$ws.Range("S16").Value = 10
$ws.Range("T16").Value = '''Västerbotten'
$ws.Range("U16").Value = '''Storuman'
$ws.Range("V16").Value = '''Lycksele lappmark'
$ws.Range("W16").Value = '''Stensele'
$ws.Range("Y16").Value = '''2023-08-08'
$ws.Range("Z16").Value = '''00:00'
$ws.Range("AA16").Value = '''2023-08-08'
$ws.Range("AB16").Value = '''00:00'
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AW16").Value = '''Isak Vahlström'
$ws.Range("AX16").Value = '''Via Isak Vahlström'

# Row 17: set to the content that was previously in row 18
$ws.Range("A17").Value = 111396318
$ws.Range("B17").Value = 96348
$ws.Range("C17").Value = '''Ovaliderad'
$ws.Range("D17").Value = '''VU'
$ws.Range("E17").Value = 220787
$ws.Range("F17").Value = '''Knärot'
$ws.Range("G17").Value = '''Goodyera repens'
$ws.Range("H17").Value = '''(L.) R. Br.'
$ws.Range("P17").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q17").Value = 625177.6865340136
$ws.Range("R17").Value = 7209552.099144561
$ws.Range("S17").Value = 10
$ws.Range("T17").Value = '''Västerbotten'
$ws.Range("U17").Value = '''Storuman'
$ws.Range("V17").Value = '''Lycksele lappmark'
$ws.Range("W17").Value = '''Stensele'
$ws.Range("Y17").Value = '''2023-08-08'
$ws.Range("Z17").Value = '''00:00'
$ws.Range("AA17").Value = '''2023-08-08'
$ws.Range("AB17").Value = '''00:00'
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AW17").Value = '''Isak Vahlström'
$ws.Range("AX17").Value = '''Via Isak Vahlström'

# Row 18: set to the content that was previously in row 9
$ws.Range("A18").Value = 111396315
$ws.Range("B18").Value = 96348
$ws.Range("C18").Value = '''Ovaliderad'
$ws.Range("D18").Value = '''VU'
$ws.Range("E18").Value = 220787
$ws.Range("F18").Value = '''Knärot'
$ws.Range("G18").Value = '''Goodyera repens'
$ws.Range("H18").Value = '''(L.) R. Br.'
$ws.Range("P18").Value = '''V Lill-Tjickuträsket, Ly lm'
$ws.Range("Q18").Value = 625167.9685939638
$ws.Range("R18").Value = 7209530.9258211
$ws.Range("S18").Value = 10
$ws.Range("T18").Value = '''Västerbotten'
$ws.Range("U18").Value = '''Storuman'
$ws.Range("V18").Value = '''Lycksele lappmark'
$ws.Range("W18").Value = '''Stensele'
$ws.Range("Y18").Value = '''2023-08-08'
$ws.Range("Z18").Value = '''00:00'
$ws.Range("AA18").Value = '''2023-08-08'
$ws.Range("AB18").Value = '''00:00'
$ws.Range("AD18").Value = $false
$ws.Range("AE18").Value = $false
$ws.Range("AG18").Value = $false
$ws.Range("AW18").Value = '''Isak Vahlström'
$ws.Range("AX18").Value = '''Via Isak Vahlström'
